$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 12 (shifts existing rows 12-20 down to 13-21,
# and fixes up the relative STDEV() formulas automatically).
$ws.Rows.Item(12).Insert()

# Populate the newly inserted row 12 with the Cunningham et al., 1999
# "predimer" data point for VEGF165:VEGFR2 (SPR method).
$ws.Range("A12").Value = "Cunningham et al., 1999"
$ws.Range("B12").Value = "VEGF165:VEGFR2"
$ws.Range("C12").Value = "VEGF165 "
$ws.Range("D12").Value = "VEGFR2 "
$ws.Range("E12").Value = "SPR"
$ws.Range("F12").Value = 4720000
$ws.Range("G12").Value = 1000000
$ws.Range("H12").Value = 0.000067
$ws.Range("I12").Value = 0.000011
$ws.Range("J12").Value = 0.0000000000145
$ws.Range("K12").Value = 0.000000000001

# Match number format of the row above used for similar kinetic-parameter rows.
$ws.Range("F12:K12").NumberFormat = "0.00E+00"

# Update the active cell selection to reflect where the edit was made.
[void]$ws.Range("K12").Select()
